$arr = New-Object 'object[,]' 21,7
$arr[0,0] = "Parameter"
$arr[0,1] = "Station"
$arr[0,2] = "Year1"
$arr[0,3] = "Year2"
$arr[0,4] = "P"
$arr[0,5] = "Estimate"
$arr[0,6] = "Change"
$arr[1,0] = "sum(precipitation_amount P1Y)"
$arr[1,1] = "SN3780"
$arr[1,2] = 1980
$arr[1,3] = 2019
$arr[1,4] = 0.138958294570753
$arr[1,5] = 3.59615384615384
$arr[1,6] = 125.865384615385
$arr[2,0] = "sum(precipitation_amount P1Y)"
$arr[2,1] = "SN18700"
$arr[2,2] = 1980
$arr[2,3] = 2019
$arr[2,4] = 0.0673677916589717
$arr[2,5] = 4.02008620689655
$arr[2,6] = 140.703017241379
$arr[3,0] = "sum(precipitation_amount P1Y)"
$arr[3,1] = "SN19710"
$arr[3,2] = 1983
$arr[3,3] = 2019
$arr[3,4] = 0.146569514589274
$arr[3,5] = 3.49074074074074
$arr[3,6] = 122.175925925926
$arr[4,0] = "sum(precipitation_amount P1Y)"
$arr[4,1] = "SN30000"
$arr[4,2] = 1980
$arr[4,3] = 2019
$arr[4,4] = 0.876850699870331
$arr[4,5] = -0.232812500000001
$arr[4,6] = -8.14843750000005
$arr[5,0] = "sum(precipitation_amount P1Y)"
$arr[5,1] = "SN30260"
$arr[5,2] = 1980
$arr[5,3] = 2015
$arr[5,4] = 0.148946024378606
$arr[5,5] = 3.80347222222222
$arr[5,6] = 133.121527777778
$arr[6,0] = "sum(precipitation_amount P1Y)"
$arr[6,1] = "SN36560"
$arr[6,2] = 1980
$arr[6,3] = 2019
$arr[6,4] = 0.0236392119138704
$arr[6,5] = 9.14074074074074
$arr[6,6] = 319.925925925926
$arr[7,0] = "sum(precipitation_amount P1Y)"
$arr[7,1] = "SN39040"
$arr[7,2] = 1980
$arr[7,3] = 2019
$arr[7,4] = 0.147464391575845
$arr[7,5] = 6.5
$arr[7,6] = 227.5
$arr[8,0] = "sum(precipitation_amount P1Y)"
$arr[8,1] = "SN43360"
$arr[8,2] = 1980
$arr[8,3] = 2017
$arr[8,4] = 0.365371548508621
$arr[8,5] = 3.48235294117646
$arr[8,6] = 121.882352941176
$arr[9,0] = "sum(precipitation_amount P1Y)"
$arr[9,1] = "SN44080"
$arr[9,2] = 1980
$arr[9,3] = 2019
$arr[9,4] = 0.711571463676106
$arr[9,5] = 1.18181818181818
$arr[9,6] = 41.3636363636364
$arr[10,0] = "sum(precipitation_amount P1Y)"
$arr[10,1] = "SN46850"
$arr[10,2] = 1980
$arr[10,3] = 2019
$arr[10,4] = 0.0994828152565819
$arr[10,5] = 12.9318181818182
$arr[10,6] = 452.613636363636
$arr[11,0] = "sum(precipitation_amount P1Y)"
$arr[11,1] = "SN51250"
$arr[11,2] = 1980
$arr[11,3] = 2019
$arr[11,4] = 0.392612310375168
$arr[11,5] = 8.42962962962963
$arr[11,6] = 295.037037037037
$arr[12,0] = "sum(precipitation_amount P1Y)"
$arr[12,1] = "SN57480"
$arr[12,2] = 1980
$arr[12,3] = 2019
$arr[12,4] = 0.506620569504461
$arr[12,5] = 4.22999999999997
$arr[12,6] = 148.049999999999
$arr[13,0] = "sum(precipitation_amount P1Y)"
$arr[13,1] = "SN63530"
$arr[13,2] = 1980
$arr[13,3] = 2019
$arr[13,4] = 0.476729058288446
$arr[13,5] = -1.76666666666667
$arr[13,6] = -61.8333333333335
$arr[14,0] = "sum(precipitation_amount P1Y)"
$arr[14,1] = "SN66210"
$arr[14,2] = 1980
$arr[14,3] = 2009
$arr[14,4] = 0.914752339084737
$arr[14,5] = 1.30000000000001
$arr[14,6] = 45.5000000000002
$arr[15,0] = "sum(precipitation_amount P1Y)"
$arr[15,1] = "SN68270"
$arr[15,2] = 1980
$arr[15,3] = 2019
$arr[15,4] = 0.345305697448102
$arr[15,5] = 2.58684210526317
$arr[15,6] = 90.539473684211
$arr[16,0] = "sum(precipitation_amount P1Y)"
$arr[16,1] = "SN78850"
$arr[16,2] = 1980
$arr[16,3] = 2007
$arr[16,4] = 0.243762192173415
$arr[16,5] = 9.04166666666667
$arr[16,6] = 316.458333333333
$arr[17,0] = "sum(precipitation_amount P1Y)"
$arr[17,1] = "SN89350"
$arr[17,2] = 1980
$arr[17,3] = 2019
$arr[17,4] = 0.301495521705862
$arr[17,5] = 1.52290969899666
$arr[17,6] = 53.301839464883
$arr[18,0] = "sum(precipitation_amount P1Y)"
$arr[18,1] = "SN93140"
$arr[18,2] = 1980
$arr[18,3] = 2017
$arr[18,4] = 0.0242392703746717
$arr[18,5] = 3.38928571428571
$arr[18,6] = 118.625
$arr[19,0] = "sum(precipitation_amount P1Y)"
$arr[19,1] = "SN96970"
$arr[19,2] = 1980
$arr[19,3] = 2018
$arr[19,4] = 0.798239812457203
$arr[19,5] = 0.445161290322581
$arr[19,6] = 15.5806451612903
$arr[20,0] = "sum(precipitation_amount P1Y)"
$arr[20,1] = "SN99500"
$arr[20,2] = 1980
$arr[20,3] = 2019
$arr[20,4] = 0.406443337928541
$arr[20,5] = 0.697142857142857
$arr[20,6] = 24.4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Precipitation")
$target = $ws.Range("A1:G21")
$target.Value = $arr
Write-Host "Wrote precipitation data"
